$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix inconsistent / misspelled data entries before removing the discs column
$ws.Range("D6").Value = "Sweden"
$ws.Range("C7").Value = "advanced"
$ws.Range("D7").Value = "Sweden"
$ws.Range("C8").Value = "pro"
$ws.Range("D8").Value = "Sweden"

# Remove the now-unused "discs" column (E) entirely
$ws.Range("E1:E8").Select() | Out-Null
$ws.Range("E1:E8").EntireColumn.Delete()
